# Update the generated answer table with the newly regenerated problems.
# Each of the 5 populated rows keeps its original 5 columns; only the
# division-problem text in each cell changes.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "91÷7=13, 0"
$t.Cell(1, 2).Range.Text = "75÷9=8, 3"
$t.Cell(1, 3).Range.Text = "60÷7=8, 4"
$t.Cell(1, 4).Range.Text = "18÷2=9, 0"
$t.Cell(1, 5).Range.Text = "34÷7=4, 6"

# Row 5
$t.Cell(5, 1).Range.Text = "60÷4=15, 0"
$t.Cell(5, 2).Range.Text = "42÷4=10, 2"
$t.Cell(5, 3).Range.Text = "65÷8=8, 1"
$t.Cell(5, 4).Range.Text = "51÷5=10, 1"
$t.Cell(5, 5).Range.Text = "44÷4=11, 0"

# Row 9
$t.Cell(9, 1).Range.Text = "80÷2=40, 0"
$t.Cell(9, 2).Range.Text = "85÷4=21, 1"
$t.Cell(9, 3).Range.Text = "28÷8=3, 4"
$t.Cell(9, 4).Range.Text = "85÷3=28, 1"
$t.Cell(9, 5).Range.Text = "59÷3=19, 2"

# Row 13
$t.Cell(13, 1).Range.Text = "60÷8=7, 4"
$t.Cell(13, 2).Range.Text = "27÷7=3, 6"
$t.Cell(13, 3).Range.Text = "96÷3=32, 0"
$t.Cell(13, 4).Range.Text = "86÷4=21, 2"
$t.Cell(13, 5).Range.Text = "23÷9=2, 5"

# Row 17
$t.Cell(17, 1).Range.Text = "77÷8=9, 5"
$t.Cell(17, 2).Range.Text = "68÷2=34, 0"
$t.Cell(17, 3).Range.Text = "85÷5=17, 0"
$t.Cell(17, 4).Range.Text = "54÷2=27, 0"
$t.Cell(17, 5).Range.Text = "12÷2=6, 0"
